$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 687, shifting existing rows 687:719 down to 688:720
$ws.Rows.Item(687).Insert()

# Populate the newly inserted row 687 with the new record
$ws.Cells.Item(687, 1).Value = 10
$ws.Cells.Item(687, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(687, 3).Value = "La Araucanía"
$ws.Cells.Item(687, 4).Value = 45041
$ws.Cells.Item(687, 5).Value = 9
$ws.Cells.Item(687, 6).Value = "Fruta"
$ws.Cells.Item(687, 7).Value = 100108
$ws.Cells.Item(687, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(687, 9).Value = 100108005
$ws.Cells.Item(687, 10).Value = "Piña"
$ws.Cells.Item(687, 11).Value = "Caramelo"
$ws.Cells.Item(687, 12).Value = "Segunda"
$ws.Cells.Item(687, 13).Value = 65
$ws.Cells.Item(687, 14).Value = 22000
$ws.Cells.Item(687, 15).Value = 22000
$ws.Cells.Item(687, 16).Value = 22000
$ws.Cells.Item(687, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(687, 18).Value = "Ecuador"
$ws.Cells.Item(687, 19).Value = 1571
$ws.Cells.Item(687, 20).Value = 14

# Ensure the date style (style index 2 / numFmtId 165) used by column D carries over correctly
$ws.Cells.Item(687, 4).NumberFormat = $ws.Cells.Item(688, 4).NumberFormat
